$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row variable holds columns B..O for one data row (rows 2..25).
$row0 = @(8.417929682717652, 4.566468218067921, 5.994567865075691, 12.49226443449704, 0, 3.651497798296613, 0, 22.96950952679996, 0, 8.643051022463588, 0, 13.73217974109775, 0, 24.76987543717602)
$row1 = @(8.149967460084351, 4.324599871906938, 5.877505468004205, 12.27668330199397, 0, 3.653589612829342, 0, 23.03508144970211, 0, 8.467371804838443, 0, 13.56968567178112, 0, 24.81144399586466)
$row2 = @(7.982500893112043, 4.168002570243659, 5.806214910329673, 12.14671911796465, 0, 3.654941667137546, 0, 23.07983653191218, 0, 8.359338476919039, 0, 13.47242771126962, 0, 24.8424351687558)
$row3 = @(7.913630780594754, 4.102182150742393, 5.777355519868758, 12.09443814507544, 0, 3.655509712370183, 0, 23.09920199934314, 0, 8.315335579210672, 0, 13.43347215954164, 0, 24.85643595483154)
$row4 = @(7.902160581173033, 4.091132530656341, 5.772576354092129, 12.08580023216239, 0, 3.655605068596373, 0, 23.1024856366137, 0, 8.308032108381356, 0, 13.42704586796351, 0, 24.85884351345593)
$row5 = @(7.981574467258847, 4.167122972210731, 5.805824865428176, 12.14601118342882, 0, 3.654949258796145, 0, 23.08009314002934, 0, 8.358744866604995, 0, 13.47189953915657, 0, 24.84261843927991)
$row6 = @(8.326210582763093, 4.484766049311015, 5.954109746328216, 12.41747913948659, 0, 3.652205044160893, 0, 22.99118455515529, 0, 8.582548629891718, 0, 13.67565780242728, 0, 24.78307152422227)
$row7 = @(8.974240102031921, 5.042524232286652, 6.247614410932792, 12.96552499852977, 0, 3.647358057346894, 0, 22.85260174687254, 0, 9.017423775184705, 0, 14.09310963333719, 0, 24.70981971236305)
$row8 = @(9.427990793556564, 5.411652284252185, 6.462393610231082, 13.37324818563349, 0, 3.644119215712593, 0, 22.77273740792197, 0, 9.33103920428209, 0, 14.40789198006435, 0, 24.68269414340063)
$row9 = @(9.628589293648011, 5.570625123270541, 6.559412444717363, 13.55890545812463, 0, 3.642714998117949, 0, 22.74120156449923, 0, 9.471756023095756, 0, 14.55224407793236, 0, 24.67617603963734)
$row10 = @(9.703644882926898, 5.629530692148534, 6.596013024432171, 13.6291607616525, 0, 3.642193143440825, 0, 22.72995152175915, 0, 9.524709374359263, 0, 14.60702240613246, 0, 24.67454623176456)
$row11 = @(9.687521685575767, 5.616901972638541, 6.588137198534491, 13.61403337465719, 0, 3.64230509504438, 0, 22.73234361053709, 0, 9.513320542257468, 0, 14.59522058682648, 0, 24.67485993327694)
$row12 = @(9.63478271031962, 5.575497289260277, 6.562426589442226, 13.56468680842846, 0, 3.642671866860953, 0, 22.74026213675791, 0, 9.476119493239514, 0, 14.55674870175763, 0, 24.67602514449153)
$row13 = @(9.602358524608405, 5.549967020516303, 6.546658937519294, 13.53445203684859, 0, 3.642897811909761, 0, 22.7452026372525, 0, 9.453287889920972, 0, 14.53319708085919, 0, 24.67684809330633)
$row14 = @(9.414757990240281, 5.401082323218058, 6.456035960223224, 13.36111221118208, 0, 3.644212371668743, 0, 22.77489504237663, 0, 9.321799227637726, 0, 14.39847712544386, 0, 24.68323738472812)
$row15 = @(9.2981291039513, 5.307449051389148, 6.400237187329886, 13.25476602850314, 0, 3.645036485265031, 0, 22.79434019644322, 0, 9.240597502511323, 0, 14.31609045332649, 0, 24.68864912649418)
$row16 = @(9.230503622426674, 5.252751983183892, 6.368080956213001, 13.19361998027015, 0, 3.645517005160366, 0, 22.80597571221091, 0, 9.193711982684293, 0, 14.26881564349588, 0, 24.69230970991635)
$row17 = @(9.207515704973037, 5.234088207429606, 6.35718396441053, 13.17292300381028, 0, 3.645680820875594, 0, 22.80999271281744, 0, 9.177808025927872, 0, 14.25282997340396, 0, 24.6936431730775)
$row18 = @(9.310601247410908, 5.317503647284925, 6.406183800000112, 13.26608509391397, 0, 3.644948083417406, 0, 22.79222351357232, 0, 9.249260632547273, 0, 14.32484945276624, 0, 24.68801632124474)
$row19 = @(9.650298557154255, 5.587694016860541, 6.569982479788046, 13.57918300824724, 0, 3.64256386902725, 0, 22.73791747725472, 0, 9.48705578716846, 0, 14.56804609928917, 0, 24.67566012996136)
$row20 = @(9.866998051572057, 5.756734623048533, 6.676211076612207, 13.7834917730795, 0, 3.641063281992074, 0, 22.70645927064641, 0, 9.640505184742825, 0, 14.72763831332609, 0, 24.67247208994657)
$row21 = @(9.751848848255028, 5.667206744907706, 6.619602794740833, 13.67450160051802, 0, 3.641858916791901, 0, 22.72287924211965, 0, 9.558802626621054, 0, 14.64241816463422, 0, 24.67372608812106)
$row22 = @(9.304964373737908, 5.312960661488801, 6.403495572512065, 13.26096775804684, 0, 3.644988028936744, 0, 22.79317904487396, 0, 9.245344655229164, 0, 14.32088922368724, 0, 24.68830070150594)
$row23 = @(8.802499637692613, 4.898722680548816, 6.16818766683616, 12.81605152346134, 0, 3.648612451923038, 0, 22.88624729283391, 0, 8.900572907807643, 0, 13.97855299283319, 0, 24.70981971236305)

$rowsData = @($row0, $row1, $row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10, $row11, $row12, $row13, $row14, $row15, $row16, $row17, $row18, $row19, $row20, $row21, $row22, $row23)

$numRows = $rowsData.Count
$numCols = 14
$arr = New-Object 'object[,]' $numRows,$numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r, $c] = $rowsData[$r][$c]
    }
}
$ws.Range("B2:O25").Value = $arr
Write-Host "Updated loading_percent values for B2:O25"
